$wb = $excel.ActiveWorkbook

# --- Create "Rebounds" sheet: copy of "Points" layout, placed right after "Assists" ---
$wb.Worksheets.Item("Points").Copy($null, $wb.Worksheets.Item("Assists"))
$wb.Worksheets.Item("Points (2)").Name = "Rebounds"
$ws = $wb.Worksheets.Item("Rebounds")
    $ws.Range("C2").Value = 1
    $ws.Range("D2").Value = 1
    $ws.Range("E2").Value = 8
    $ws.Range("F2").Value = 0
    $ws.Range("G2").Value = 2
    $ws.Range("H2").Value = 3
    $ws.Range("I2").Value = 6
    $ws.Range("J2").Value = 0
    $ws.Range("K2").Value = 1
    $ws.Range("L2").Value = 6
    $ws.Range("M2").Value = 4
    $ws.Range("N2").Value = 0
    $ws.Range("O2").Value = 1
    $ws.Range("P2").Value = 4
    $ws.Range("C3").Value = 0
    $ws.Range("D3").Value = 1
    $ws.Range("E3").Value = 0
    $ws.Range("F3").Value = 0
    $ws.Range("G3").Value = 2
    $ws.Range("H3").Value = 2
    $ws.Range("I3").Value = 3
    $ws.Range("J3").Value = 12
    $ws.Range("K3").Value = 5
    $ws.Range("L3").Value = 1
    $ws.Range("M3").Value = 2
    $ws.Range("N3").Value = 0
    $ws.Range("O3").Value = 3
    $ws.Range("P3").Value = 0
    $ws.Range("C4").Value = 3
    $ws.Range("D4").Value = 1
    $ws.Range("E4").Value = 0
    $ws.Range("F4").Value = 0
    $ws.Range("G4").Value = 4
    $ws.Range("H4").Value = 0
    $ws.Range("I4").Value = 6
    $ws.Range("J4").Value = 14
    $ws.Range("K4").Value = 2
    $ws.Range("L4").Value = 6
    $ws.Range("M4").Value = 2
    $ws.Range("N4").Value = 0
    $ws.Range("O4").Value = 2
    $ws.Range("P4").Value = 0
    $ws.Range("C5").Value = 0
    $ws.Range("D5").Value = 2
    $ws.Range("E5").Value = 0
    $ws.Range("F5").Value = 0
    $ws.Range("G5").Value = 4
    $ws.Range("H5").Value = 0
    $ws.Range("I5").Value = 1
    $ws.Range("J5").Value = 18
    $ws.Range("K5").Value = 5
    $ws.Range("L5").Value = 9
    $ws.Range("M5").Value = 4
    $ws.Range("N5").Value = 1
    $ws.Range("O5").Value = 2
    $ws.Range("P5").Value = 0
    $ws.Range("C6").Value = 0
    $ws.Range("D6").Value = 1
    $ws.Range("E6").Value = 0
    $ws.Range("F6").Value = 0
    $ws.Range("G6").Value = 3
    $ws.Range("H6").Value = 0
    $ws.Range("I6").Value = 6
    $ws.Range("J6").Value = 11
    $ws.Range("K6").Value = 2
    $ws.Range("L6").Value = 2
    $ws.Range("M6").Value = 3
    $ws.Range("N6").Value = 2
    $ws.Range("O6").Value = 1
    $ws.Range("P6").Value = 0
    $ws.Range("C7").Value = 1
    $ws.Range("D7").Value = 0
    $ws.Range("E7").Value = 3
    $ws.Range("F7").Value = 0
    $ws.Range("G7").Value = 4
    $ws.Range("H7").Value = 0
    $ws.Range("I7").Value = 4
    $ws.Range("J7").Value = 13
    $ws.Range("K7").Value = 5
    $ws.Range("L7").Value = 5
    $ws.Range("M7").Value = 1
    $ws.Range("N7").Value = 0
    $ws.Range("O7").Value = 1
    $ws.Range("P7").Value = 1
    $ws.Range("C8").Value = 0
    $ws.Range("D8").Value = 0
    $ws.Range("E8").Value = 1
    $ws.Range("F8").Value = 0
    $ws.Range("G8").Value = 5
    $ws.Range("H8").Value = 0
    $ws.Range("I8").Value = 3
    $ws.Range("J8").Value = 17
    $ws.Range("K8").Value = 4
    $ws.Range("L8").Value = 12
    $ws.Range("M8").Value = 2
    $ws.Range("N8").Value = 0
    $ws.Range("O8").Value = 3
    $ws.Range("P8").Value = 0
    $ws.Range("C9").Value = 0
    $ws.Range("D9").Value = 4
    $ws.Range("E9").Value = 3
    $ws.Range("F9").Value = 3
    $ws.Range("G9").Value = 0
    $ws.Range("H9").Value = 0
    $ws.Range("I9").Value = 2
    $ws.Range("J9").Value = 0
    $ws.Range("K9").Value = 7
    $ws.Range("L9").Value = 16
    $ws.Range("M9").Value = 6
    $ws.Range("N9").Value = 1
    $ws.Range("O9").Value = 2
    $ws.Range("P9").Value = 3
    $ws.Range("C10").Value = 0
    $ws.Range("D10").Value = 3
    $ws.Range("E10").Value = 3
    $ws.Range("F10").Value = 4
    $ws.Range("G10").Value = 2
    $ws.Range("H10").Value = 0
    $ws.Range("I10").Value = 4
    $ws.Range("J10").Value = 0
    $ws.Range("K10").Value = 2
    $ws.Range("L10").Value = 7
    $ws.Range("M10").Value = 1
    $ws.Range("N10").Value = 1
    $ws.Range("O10").Value = 1
    $ws.Range("P10").Value = 10
    $ws.Range("C11").Value = 0
    $ws.Range("D11").Value = 0
    $ws.Range("E11").Value = 5
    $ws.Range("F11").Value = 10
    $ws.Range("G11").Value = 4
    $ws.Range("H11").Value = 0
    $ws.Range("I11").Value = 3
    $ws.Range("J11").Value = 13
    $ws.Range("K11").Value = 3
    $ws.Range("L11").Value = 0
    $ws.Range("M11").Value = 1
    $ws.Range("N11").Value = 1
    $ws.Range("O11").Value = 3
    $ws.Range("P11").Value = 0

# --- Create "3PM" sheet: copy of "Points" layout, placed right after "Rebounds" ---
$wb.Worksheets.Item("Points").Copy($null, $wb.Worksheets.Item("Rebounds"))
$wb.Worksheets.Item("Points (2)").Name = "3PM"
$ws = $wb.Worksheets.Item("3PM")
    $ws.Range("C2").Value = 0
    $ws.Range("D2").Value = 3
    $ws.Range("E2").Value = 0
    $ws.Range("F2").Value = 0
    $ws.Range("G2").Value = 2
    $ws.Range("H2").Value = 0
    $ws.Range("I2").Value = 1
    $ws.Range("J2").Value = 0
    $ws.Range("K2").Value = 2
    $ws.Range("L2").Value = 0
    $ws.Range("M2").Value = 0
    $ws.Range("N2").Value = 0
    $ws.Range("O2").Value = 1
    $ws.Range("P2").Value = 0
    $ws.Range("C3").Value = 0
    $ws.Range("D3").Value = 6
    $ws.Range("E3").Value = 0
    $ws.Range("F3").Value = 0
    $ws.Range("G3").Value = 4
    $ws.Range("H3").Value = 1
    $ws.Range("I3").Value = 0
    $ws.Range("J3").Value = 0
    $ws.Range("K3").Value = 1
    $ws.Range("L3").Value = 1
    $ws.Range("M3").Value = 0
    $ws.Range("N3").Value = 0
    $ws.Range("O3").Value = 1
    $ws.Range("P3").Value = 0
    $ws.Range("C4").Value = 0
    $ws.Range("D4").Value = 3
    $ws.Range("E4").Value = 0
    $ws.Range("F4").Value = 0
    $ws.Range("G4").Value = 4
    $ws.Range("H4").Value = 0
    $ws.Range("I4").Value = 2
    $ws.Range("J4").Value = 1
    $ws.Range("K4").Value = 3
    $ws.Range("L4").Value = 4
    $ws.Range("M4").Value = 0
    $ws.Range("N4").Value = 0
    $ws.Range("O4").Value = 1
    $ws.Range("P4").Value = 0
    $ws.Range("C5").Value = 0
    $ws.Range("D5").Value = 0
    $ws.Range("E5").Value = 0
    $ws.Range("F5").Value = 0
    $ws.Range("G5").Value = 4
    $ws.Range("H5").Value = 0
    $ws.Range("I5").Value = 0
    $ws.Range("J5").Value = 0
    $ws.Range("K5").Value = 0
    $ws.Range("L5").Value = 4
    $ws.Range("M5").Value = 0
    $ws.Range("N5").Value = 0
    $ws.Range("O5").Value = 2
    $ws.Range("P5").Value = 0
    $ws.Range("C6").Value = 0
    $ws.Range("D6").Value = 1
    $ws.Range("E6").Value = 0
    $ws.Range("F6").Value = 0
    $ws.Range("G6").Value = 2
    $ws.Range("H6").Value = 0
    $ws.Range("I6").Value = 0
    $ws.Range("J6").Value = 1
    $ws.Range("K6").Value = 1
    $ws.Range("L6").Value = 1
    $ws.Range("M6").Value = 0
    $ws.Range("N6").Value = 0
    $ws.Range("O6").Value = 1
    $ws.Range("P6").Value = 0
    $ws.Range("C7").Value = 0
    $ws.Range("D7").Value = 0
    $ws.Range("E7").Value = 0
    $ws.Range("F7").Value = 0
    $ws.Range("G7").Value = 6
    $ws.Range("H7").Value = 0
    $ws.Range("I7").Value = 0
    $ws.Range("J7").Value = 0
    $ws.Range("K7").Value = 3
    $ws.Range("L7").Value = 0
    $ws.Range("M7").Value = 0
    $ws.Range("N7").Value = 0
    $ws.Range("O7").Value = 1
    $ws.Range("P7").Value = 0
    $ws.Range("C8").Value = 0
    $ws.Range("D8").Value = 0
    $ws.Range("E8").Value = 0
    $ws.Range("F8").Value = 0
    $ws.Range("G8").Value = 1
    $ws.Range("H8").Value = 0
    $ws.Range("I8").Value = 0
    $ws.Range("J8").Value = 0
    $ws.Range("K8").Value = 0
    $ws.Range("L8").Value = 3
    $ws.Range("M8").Value = 0
    $ws.Range("N8").Value = 0
    $ws.Range("O8").Value = 5
    $ws.Range("P8").Value = 0
    $ws.Range("C9").Value = 0
    $ws.Range("D9").Value = 2
    $ws.Range("E9").Value = 3
    $ws.Range("F9").Value = 0
    $ws.Range("G9").Value = 0
    $ws.Range("H9").Value = 0
    $ws.Range("I9").Value = 0
    $ws.Range("J9").Value = 0
    $ws.Range("K9").Value = 3
    $ws.Range("L9").Value = 3
    $ws.Range("M9").Value = 0
    $ws.Range("N9").Value = 0
    $ws.Range("O9").Value = 0
    $ws.Range("P9").Value = 0
    $ws.Range("C10").Value = 0
    $ws.Range("D10").Value = 3
    $ws.Range("E10").Value = 0
    $ws.Range("F10").Value = 0
    $ws.Range("G10").Value = 2
    $ws.Range("H10").Value = 0
    $ws.Range("I10").Value = 4
    $ws.Range("J10").Value = 0
    $ws.Range("K10").Value = 0
    $ws.Range("L10").Value = 5
    $ws.Range("M10").Value = 0
    $ws.Range("N10").Value = 0
    $ws.Range("O10").Value = 1
    $ws.Range("P10").Value = 0
    $ws.Range("C11").Value = 0
    $ws.Range("D11").Value = 0
    $ws.Range("E11").Value = 1
    $ws.Range("F11").Value = 0
    $ws.Range("G11").Value = 6
    $ws.Range("H11").Value = 0
    $ws.Range("I11").Value = 3
    $ws.Range("J11").Value = 0
    $ws.Range("K11").Value = 0
    $ws.Range("L11").Value = 2
    $ws.Range("M11").Value = 0
    $ws.Range("N11").Value = 0
    $ws.Range("O11").Value = 1
    $ws.Range("P11").Value = 0

# --- Create "Avg Rebounds" sheet: copy of "Avg Points" layout, placed right after "Avg Assists" ---
$wb.Worksheets.Item("Avg Points").Copy($null, $wb.Worksheets.Item("Avg Assists"))
$wb.Worksheets.Item("Avg Points (2)").Name = "Avg Rebounds"
$ws = $wb.Worksheets.Item("Avg Rebounds")
$ws.Range("B1").Value = "Avg Rebounds"
    $ws.Range("A2").Value = "Domantas Sabonis"
    $ws.Range("B2").Value = 14
    $ws.Range("A3").Value = "Russell Westbrook"
    $ws.Range("B3").Value = 6.4
    $ws.Range("A4").Value = "Precious Achiuwa"
    $ws.Range("B4").Value = 5.666666666666667
    $ws.Range("A5").Value = "Nique Clifford"
    $ws.Range("B5").Value = 3.833333333333333
    $ws.Range("A6").Value = "DeMar DeRozan"
    $ws.Range("B6").Value = 3.8
    $ws.Range("A7").Value = "Dennis Schröder"
    $ws.Range("B7").Value = 3.6
    $ws.Range("A8").Value = "Maxime Raynaud"
    $ws.Range("B8").Value = 3.6
    $ws.Range("A9").Value = "Zach LaVine"
    $ws.Range("B9").Value = 3.333333333333333
    $ws.Range("A10").Value = "Drew Eubanks"
    $ws.Range("B10").Value = 2.6
    $ws.Range("A11").Value = "Keon Ellis"
    $ws.Range("B11").Value = 1.9
    $ws.Range("A12").Value = "Dylan Cardwell"
    $ws.Range("B12").Value = 1.666666666666667
    $ws.Range("A13").Value = "Dario Šarić"
    $ws.Range("B13").Value = 1.666666666666667
    $ws.Range("A14").Value = "Malik Monk"
    $ws.Range("B14").Value = 1.625
    $ws.Range("A15").Value = "Devin Carter"
    $ws.Range("B15").Value = 0.8571428571428571

# --- Create "Avg 3PM" sheet: copy of "Avg Points" layout, placed right after "Avg Rebounds" ---
$wb.Worksheets.Item("Avg Points").Copy($null, $wb.Worksheets.Item("Avg Rebounds"))
$wb.Worksheets.Item("Avg Points (2)").Name = "Avg 3PM"
$ws = $wb.Worksheets.Item("Avg 3PM")
$ws.Range("B1").Value = "Avg 3PM"
    $ws.Range("A2").Value = "Zach LaVine"
    $ws.Range("B2").Value = 3.444444444444445
    $ws.Range("A3").Value = "Russell Westbrook"
    $ws.Range("B3").Value = 2.3
    $ws.Range("A4").Value = "Malik Monk"
    $ws.Range("B4").Value = 2.25
    $ws.Range("A5").Value = "Keon Ellis"
    $ws.Range("B5").Value = 1.4
    $ws.Range("A6").Value = "Dennis Schröder"
    $ws.Range("B6").Value = 1.3
    $ws.Range("A7").Value = "DeMar DeRozan"
    $ws.Range("B7").Value = 1
    $ws.Range("A8").Value = "Nique Clifford"
    $ws.Range("B8").Value = 0.6666666666666666
    $ws.Range("A9").Value = "Dario Šarić"
    $ws.Range("B9").Value = 0.3333333333333333
    $ws.Range("A10").Value = "Domantas Sabonis"
    $ws.Range("B10").Value = 0.2857142857142857
    $ws.Range("A11").Value = "Dylan Cardwell"
    $ws.Range("B11").Value = 0
    $ws.Range("A12").Value = "Precious Achiuwa"
    $ws.Range("B12").Value = 0
    $ws.Range("A13").Value = "Drew Eubanks"
    $ws.Range("B13").Value = 0
    $ws.Range("A14").Value = "Devin Carter"
    $ws.Range("B14").Value = 0
    $ws.Range("A15").Value = "Maxime Raynaud"
    $ws.Range("B15").Value = 0
